$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue 2 4 "66.841.70"
Set-TextValue 2 5 "  +2.77%  "

# Row 3 - Ethereum
Set-TextValue 3 4 "3.701.28"
Set-TextValue 3 5 "  +5.23%  "

# Row 4 - TetherUSD
Set-TextValue 4 5 "  +0.06%  "

# Row 5 - BNB
Set-TextValue 5 4 "419.94"
Set-TextValue 5 5 "  -0.43%  "

# Row 6 - Solana
Set-TextValue 6 4 "130.64"
Set-TextValue 6 5 "  -0.98%  "

# Row 7 - LidoStakedEther
Set-TextValue 7 4 "3.692.18"
Set-TextValue 7 5 "  +5.15%  "

# Row 8 - XRP
Set-TextValue 8 5 "  +1.21%  "

# Row 9 - USDC
Set-TextValue 9 5 "  +0.04%  "

# Row 10 - Cardano
Set-TextValue 10 4 "0.769"
Set-TextValue 10 5 "  -2.52%  "

# Row 11 - Dogecoin
Set-TextValue 11 5 "  +10.89%  "

# Row 12 - ShibaInu
Set-TextValue 12 4 "0.0000400"
Set-TextValue 12 5 "  +48.75%  "

# Row 13 - Avalanche
Set-TextValue 13 4 "43.20"
Set-TextValue 13 5 "  +0.24%  "

# Row 14 - Polkadot
Set-TextValue 14 4 "10.72"
Set-TextValue 14 5 "  +6.68%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue 15 4 "4.277.05"
Set-TextValue 15 5 "  +4.90%  "

# Row 16 - TRON
Set-TextValue 16 5 "  -0.85%  "

# Row 17 - Chainlink
Set-TextValue 17 4 "20.64"
Set-TextValue 17 5 "  -0.21%  "

# Row 18 - WrappedEther
Set-TextValue 18 4 "3.704.73"
Set-TextValue 18 5 "  +5.38%  "

# Row 19 - Uniswap
Set-TextValue 19 4 "13.19"
Set-TextValue 19 5 "  +5.15%  "

# Row 20 - Polygon
Set-TextValue 20 5 "  +1.99%  "

# Row 21 - WrappedBTC
Set-TextValue 21 4 "66.831.49"
Set-TextValue 21 5 "  +3.01%  "

# Row 22 - BitcoinCash
Set-TextValue 22 4 "445.44"
Set-TextValue 22 5 "  -2.12%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextValue 23 4 "16.48"
Set-TextValue 23 5 "  +23.52%  "

# Row 24 - Litecoin
Set-TextValue 24 4 "90.36"
Set-TextValue 24 5 "  -1.23%  "

# Row 25 - ImmutableX
Set-TextValue 25 5 "  -2.28%  "

# Row 26 - EthereumClassic
Set-TextValue 26 5 "  +9.45%  "

# Row 27 - Filecoin
Set-TextValue 27 4 "10.30"
Set-TextValue 27 5 "  +0.73%  "

# Row 28 - PancakeSwap
Set-TextValue 28 4 "3.33"
Set-TextValue 28 5 "  -1.15%  "

# Row 29 - LEO
Set-TextValue 29 4 "5.06"
Set-TextValue 29 5 "  +5.29%  "

# Row 30 - Hedera
Set-TextValue 30 5 "  +10.86%  "

# Row 31 - Cosmos
Set-TextValue 31 5 "  +1.80%  "

# Row 32 - Toncoin
Set-TextValue 32 5 "  +2.38%  "

# Row 34 - Kaspa
Set-TextValue 34 5 "  +0.84%  "

# Row 35 - InjectiveProtocol
Set-TextValue 35 4 "41.53"
Set-TextValue 35 5 "  +3.37%  "

# Row 36 - OKB
Set-TextValue 36 4 "57.22"
Set-TextValue 36 5 "  -0.61%  "

# Row 37 - Dai
Set-TextValue 37 5 "  -0.10%  "

# Row 38 - VeChain
Set-TextValue 38 4 "0.0496"
Set-TextValue 38 5 "  -2.63%  "

# Row 39 - PEPE
Set-TextValue 39 4 "0.0₃0745"
Set-TextValue 39 5 "  +2.33%  "

# Row 40 - now ThetaToken (was Stellar)
Set-TextValue 40 2 "ThetaToken"
Set-TextValue 40 3 "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue 40 4 "3.09"
Set-TextValue 40 5 "  +33.23%  "

# Row 41 - now Stellar (was ThetaToken)
Set-TextValue 41 2 "Stellar"
Set-TextValue 41 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 41 4 "0.152"
Set-TextValue 41 5 "  +4.47%  "

# Row 42 - EnergySwap
Set-TextValue 42 4 "29.73"
Set-TextValue 42 5 "  +34.40%  "

# Row 43 - FirstDigitalUSD
Set-TextValue 43 4 "0.997"
Set-TextValue 43 5 "  -0.19%  "

# Row 44 - LidoDAOToken
Set-TextValue 44 5 "  +2.65%  "

# Row 45 - Monero
Set-TextValue 45 4 "149.26"
Set-TextValue 45 5 "  +2.04%  "

# Row 46 - ARBITRUM
Set-TextValue 46 5 "  +4.70%  "

# Row 47 - Stacks
Set-TextValue 47 4 "2.90"
Set-TextValue 47 5 "  -6.32%  "

# Row 48 - WEMIXToken
Set-TextValue 48 4 "2.66"
Set-TextValue 48 5 "  -4.91%  "

# Row 49 - NEARProtocol
Set-TextValue 49 4 "4.35"
Set-TextValue 49 5 "  -4.84%  "

# Row 50 - TheGraph
Set-TextValue 50 4 "0.307"
Set-TextValue 50 5 "  -3.32%  "

# Row 51 - Cronos
Set-TextValue 51 4 "0.161"
Set-TextValue 51 5 "  +12.41%  "
